# Update the "Price" (D) and "Volume(1h)" (E) columns of the cryptos list with the
# refreshed figures from the latest GitHub Actions run.
#
# Note: some Price values (e.g. "380.93") look like plain decimal numbers. A bare
# string assignment would be auto-coerced by Excel into a floating point number
# (losing the exact text/trailing zeros and introducing binary float noise such as
# 380.93000000000001). To keep these as literal text -- matching how they were
# authored in the workbook -- we prefix them with a leading apostrophe, which is
# the standard Excel mechanism for forcing text storage of a number-looking value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.712.77"
$ws.Range("E2").Value = "  +1.95%  "

$ws.Range("D3").Value = "3.036.17"
$ws.Range("E3").Value = "  +4.02%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'380.93"
$ws.Range("E5").Value = "  +1.70%  "

$ws.Range("D6").Value = "'103.54"
$ws.Range("E6").Value = "  +4.06%  "

$ws.Range("E7").Value = "  +2.27%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  +3.99%  "

$ws.Range("D10").Value = "'36.92"
$ws.Range("E10").Value = "  +3.76%  "

$ws.Range("E11").Value = "  -0.03%  "

$ws.Range("D12").Value = "'0.0860"
$ws.Range("E12").Value = "  +1.91%  "

$ws.Range("D13").Value = "3.505.21"
$ws.Range("E13").Value = "  +3.70%  "

$ws.Range("D14").Value = "'18.60"
$ws.Range("E14").Value = "  +3.63%  "

$ws.Range("D15").Value = "'7.80"
$ws.Range("E15").Value = "  +2.48%  "

$ws.Range("D16").Value = "3.060.21"
$ws.Range("E16").Value = "  +4.18%  "

$ws.Range("D17").Value = "'0.999"
$ws.Range("E17").Value = "  +0.71%  "

$ws.Range("D18").Value = "'10.84"
$ws.Range("E18").Value = "  -9.41%  "

$ws.Range("D19").Value = "51.796.31"
$ws.Range("E19").Value = "  +2.21%  "

$ws.Range("E20").Value = "  +2.98%  "

$ws.Range("D21").Value = "'12.56"
$ws.Range("E21").Value = "  +2.68%  "

$ws.Range("E22").Value = "  +2.68%  "

$ws.Range("D23").Value = "'70.40"
$ws.Range("E23").Value = "  +1.69%  "

$ws.Range("D24").Value = "'269.02"
$ws.Range("E24").Value = "  +1.36%  "

$ws.Range("D25").Value = "'3.19"
$ws.Range("E25").Value = "  +2.01%  "

$ws.Range("D26").Value = "'8.30"
$ws.Range("E26").Value = "  +6.40%  "

$ws.Range("D27").Value = "'7.56"
$ws.Range("E27").Value = "  +6.53%  "

$ws.Range("E28").Value = "  +6.40%  "

$ws.Range("D29").Value = "'26.31"
$ws.Range("E29").Value = "  +3.94%  "

$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("E31").Value = "  +2.07%  "

$ws.Range("E32").Value = "  +4.75%  "

$ws.Range("D33").Value = "'34.44"
$ws.Range("E33").Value = "  +4.14%  "

$ws.Range("D34").Value = "'51.25"
$ws.Range("E34").Value = "  +1.96%  "

$ws.Range("E35").Value = "  +6.64%  "

$ws.Range("E36").Value = "  +0.66%  "

$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("E38").Value = "  +9.01%  "

$ws.Range("D39").Value = "'17.22"
$ws.Range("E39").Value = "  +5.97%  "

$ws.Range("D40").Value = "'2.63"
$ws.Range("E40").Value = "  +9.72%  "

$ws.Range("D41").Value = "'0.284"
$ws.Range("E41").Value = "  +10.81%  "

$ws.Range("E42").Value = "  +5.02%  "

$ws.Range("D43").Value = "'0.116"
$ws.Range("E43").Value = "  +1.58%  "

$ws.Range("D44").Value = "'127.01"
$ws.Range("E44").Value = "  +6.39%  "

$ws.Range("E45").Value = "  +14.26%  "

$ws.Range("D46").Value = "'21.95"
$ws.Range("E46").Value = "  +6.07%  "

$ws.Range("E47").Value = "  +0.64%  "

$ws.Range("D48").Value = "'2.40"
$ws.Range("E48").Value = "  +2.76%  "

$ws.Range("D49").Value = "2.043.39"
$ws.Range("E49").Value = "  +2.85%  "

$ws.Range("D50").Value = "3.337.87"
$ws.Range("E50").Value = "  +3.99%  "

$ws.Range("D51").Value = "'0.0324"
$ws.Range("E51").Value = "  +3.81%  "
